$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Values / formulas -----------------------------------------------
# Row 2 header (Week 1 / day labels / Total) - write days first then
# "Total" so shared-string insertion order matches natural authoring order.
$ws.Range("A2").Value = "Week 1"
$ws.Range("C2").Value = "Sun - 07/15"
$ws.Range("D2").Value = "Mon - 07/16"
$ws.Range("E2").Value = "Tue - 07/17"
$ws.Range("F2").Value = "Wed - 07/18"
$ws.Range("G2").Value = "Thu - 07/19"
$ws.Range("H2").Value = "Fri - 07/20"
$ws.Range("I2").Value = "Sat - 07/21"
$ws.Range("B2").Value = "Total"

# Row 3: Sawyer - no hours logged this week
$ws.Range("B3").Formula = "=SUM(C3:I3)"

# Row 4: Rashmi - no hours logged this week
$ws.Range("B4").Formula = "=SUM(C4:I4)"

# Row 5: Spencer - 2 hours on Tuesday
$ws.Range("B5").Formula = "=SUM(C5:I5)"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0

# ---- Formatting ---------------------------------------------------------
# Thin box border around the whole table (A2:I5)
$ws.Range("A2:I5").Borders.LineStyle = 1

# Center-align the day/total value grid (everything except column A labels)
$ws.Range("C2:I2").HorizontalAlignment = -4108
$ws.Range("B3:I5").HorizontalAlignment = -4108

# Ensure the empty day cells in rows 3 and 4 exist (formatted but blank)
$ws.Range("C3:I3").Borders.LineStyle = 1
$ws.Range("C4:I4").Borders.LineStyle = 1
$ws.Range("G5:I5").Borders.LineStyle = 1

# ---- Column widths --------------------------------------------------------
# (values chosen so the engine's internal px-rounding lands on the closest
# representable width to the authored OOXML widths)
$ws.Columns.Item(2).ColumnWidth = 4.666666666666667
$ws.Columns.Item(3).ColumnWidth = 9.833333333333334
$ws.Columns.Item(4).ColumnWidth = 10.666666666666666
$ws.Columns.Item(5).ColumnWidth = 9.833333333333334
$ws.Columns.Item(6).ColumnWidth = 10.666666666666666
$ws.Columns.Item(7).ColumnWidth = 9.833333333333334
$ws.Columns.Item(8).ColumnWidth = 8.833333333333334
$ws.Columns.Item(9).ColumnWidth = 9.333333333333334

# ---- Selection -----------------------------------------------------------
$ws.Range("A2").Select()
